$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) column values, forcing text storage to preserve exact formatting
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.417.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.453.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "528.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.556"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.460.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0992"
$ws.Range("D10").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.889.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.298.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.462.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "320.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.47"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.408"
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0753"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.25"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.35"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.54"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.807"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.54"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "274.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.82"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.586"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0923"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "120.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.12"
$ws.Range("D51").Style = "Normal"

# Update remaining text/percentage columns (B, C, E)
$ws.Range("E2").Value = "  -2.83%  "
$ws.Range("E3").Value = "  -3.73%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -2.62%  "
$ws.Range("E9").Value = "  -4.24%  "
$ws.Range("E10").Value = "  -1.94%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("E12").Value = "  -2.57%  "
$ws.Range("E13").Value = "  -5.49%  "
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("E16").Value = "  -5.92%  "
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("E18").Value = "  -3.92%  "
$ws.Range("E19").Value = "  -4.89%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("E25").Value = "  -6.01%  "
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -1.08%  "
$ws.Range("E28").Value = "  -7.67%  "
$ws.Range("E29").Value = "  -5.35%  "
$ws.Range("E30").Value = "  -7.97%  "
$ws.Range("E31").Value = "  -3.73%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -7.22%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("E35").Value = "  -8.54%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E36").Value = "  -2.88%  "
$ws.Range("E37").Value = "  -8.57%  "
$ws.Range("E38").Value = "  -5.83%  "
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E42").Value = "  -9.20%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("E43").Value = "  -9.59%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("E47").Value = "  -1.65%  "
$ws.Range("E48").Value = "  -5.16%  "
$ws.Range("E49").Value = "  -2.71%  "
$ws.Range("E50").Value = "  -4.60%  "
$ws.Range("E51").Value = "  -6.18%  "
